# Auto-generated edit script: update FFXIV leve profit values
# per scheduled runner refresh (see commit message).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2188.75
$ws.Range("I15").Value = 2188.75
$ws.Range("K15").Value = 6566.25
$ws.Range("M15").Value = -6397.25
$ws.Range("H17").Value = 1860.091
$ws.Range("J17").Value = 1860.091
$ws.Range("L17").Value = 5580.272999999999
$ws.Range("N17").Value = -5916.272999999999
$ws.Range("H112").Value = 2402.5293
$ws.Range("I112").Value = 1795
$ws.Range("J112").Value = 2483.5334
$ws.Range("K112").Value = 5385
$ws.Range("L112").Value = 7450.600199999999
$ws.Range("M112").Value = -4277
$ws.Range("N112").Value = -9666.600199999999
$ws.Range("H132").Value = 2823.7
$ws.Range("I132").Value = 2104.75
$ws.Range("J132").Value = 5699.5
$ws.Range("K132").Value = 6314.25
$ws.Range("L132").Value = 17098.5
$ws.Range("M132").Value = -3784.25
$ws.Range("N132").Value = -22158.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 39800
$ws.Range("J37").Value = 39800
$ws.Range("L37").Value = 39800
$ws.Range("N37").Value = -40346
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = $null
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4788
$ws.Range("H74").Value = 8094.5
$ws.Range("I74").Value = 8094.5
$ws.Range("K74").Value = 8094.5
$ws.Range("M74").Value = -7220.5
$ws.Range("H77").Value = 8094.5
$ws.Range("I77").Value = 8094.5
$ws.Range("K77").Value = 40472.5
$ws.Range("M77").Value = -36104.5
$ws.Range("H102").Value = 2581.5557
$ws.Range("I102").Value = 2526.9285
$ws.Range("K102").Value = 2526.9285
$ws.Range("M102").Value = -904.9285
$ws.Range("H132").Value = 3993.6667
$ws.Range("I132").Value = 4237
$ws.Range("J132").Value = 3799
$ws.Range("K132").Value = 12711
$ws.Range("L132").Value = 11397
$ws.Range("M132").Value = -10181
$ws.Range("N132").Value = -16457
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1747.5
$ws.Range("I99").Value = 2495
$ws.Range("K99").Value = 2495
$ws.Range("M99").Value = -997
$ws.Range("H134").Value = 4987.8335
$ws.Range("I134").Value = 4987.8335
$ws.Range("K134").Value = 14963.5005
$ws.Range("M134").Value = -12428.5005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1872.5
$ws.Range("I31").Value = 1568.7142
$ws.Range("K31").Value = 1568.7142
$ws.Range("M31").Value = -1273.7142
$ws.Range("H34").Value = 1872.5
$ws.Range("I34").Value = 1568.7142
$ws.Range("K34").Value = 1568.7142
$ws.Range("M34").Value = -1366.7142
$ws.Range("H35").Value = 15210
$ws.Range("I35").Value = 920
$ws.Range("J35").Value = 29500
$ws.Range("K35").Value = 920
$ws.Range("L35").Value = 29500
$ws.Range("M35").Value = -626
$ws.Range("N35").Value = -30088
$ws.Range("H58").Value = 3419.2856
$ws.Range("I58").Value = 3459.2307
$ws.Range("K58").Value = 3459.2307
$ws.Range("M58").Value = -3256.2307
$ws.Range("H132").Value = 3294.4
$ws.Range("I132").Value = 2774.3333
$ws.Range("J132").Value = 4074.5
$ws.Range("K132").Value = 8322.999899999999
$ws.Range("L132").Value = 12223.5
$ws.Range("M132").Value = -5792.999899999999
$ws.Range("N132").Value = -17283.5
$ws.Range("H136").Value = 3419.2856
$ws.Range("I136").Value = 3459.2307
$ws.Range("K136").Value = 10377.6921
$ws.Range("M136").Value = -7827.6921

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1800.7142
$ws.Range("I5").Value = 1767.5
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 5302.5
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = -5190.5
$ws.Range("N5").Value = -6224
$ws.Range("H32").Value = 178.33333
$ws.Range("I32").Value = 17.5
$ws.Range("K32").Value = 52.5
$ws.Range("M32").Value = 230.5
$ws.Range("H33").Value = 108
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = $null
$ws.Range("H39").Value = 15625
$ws.Range("J39").Value = 15625
$ws.Range("L39").Value = 46875
$ws.Range("N39").Value = -47463
$ws.Range("H50").Value = 100000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 100000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 300000
$ws.Range("M50").Value = $null
$ws.Range("N50").Value = -300962
$ws.Range("H53").Value = 100000
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 100000
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 300000
$ws.Range("M53").Value = $null
$ws.Range("N53").Value = -300962
$ws.Range("H61").Value = 4
$ws.Range("I61").Value = 4
$ws.Range("K61").Value = 12
$ws.Range("M61").Value = 203
$ws.Range("H104").Value = 5048.9614
$ws.Range("J104").Value = 5090.92
$ws.Range("L104").Value = 15272.76
$ws.Range("N104").Value = -20514.76
$ws.Range("H135").Value = 1800.7142
$ws.Range("I135").Value = 1767.5
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 15907.5
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -13372.5
$ws.Range("N135").Value = -23070

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = $null
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = $null
$ws.Range("H102").Value = 2644.4285
$ws.Range("I102").Value = 2222.2
$ws.Range("J102").Value = 3700
$ws.Range("K102").Value = 2222.2
$ws.Range("L102").Value = 3700
$ws.Range("M102").Value = -600.1999999999998
$ws.Range("N102").Value = -6944
$ws.Range("H122").Value = 2671.9092
$ws.Range("I122").Value = 2654.6667
$ws.Range("K122").Value = 7964.000100000001
$ws.Range("M122").Value = -5514.000100000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1131.6666
$ws.Range("I107").Value = 1158
$ws.Range("K107").Value = 3474
$ws.Range("M107").Value = -1554
$ws.Range("H122").Value = 1587.125
$ws.Range("I122").Value = 1587.125
$ws.Range("K122").Value = 4761.375
$ws.Range("M122").Value = -2311.375
$ws.Range("H132").Value = 2127.2727
$ws.Range("I132").Value = 1489.1111
$ws.Range("K132").Value = 4467.3333
$ws.Range("M132").Value = -1937.3333
$ws.Range("H136").Value = 6680
$ws.Range("I136").Value = 7850
$ws.Range("K136").Value = 23550
$ws.Range("M136").Value = -21000
